$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set string values in the exact order needed to reproduce shared-string indices:
# 2: Sediment subduction (B1)
# 3: sed (A3)
# 4: Sample sediment grid (C1)
# 5: syn (A4)
# 6: Reconstructed motions (D1)
# 7: sediment thickness (C3)
$ws.Range("B1").Value = "Sediment subduction"
$ws.Range("A3").Value = "sed"
$ws.Range("C1").Value = "Sample sediment grid"
$ws.Range("A4").Value = "syn"
$ws.Range("D1").Value = "Reconstructed motions"
$ws.Range("C3").Value = "sediment thickness"

# Booleans
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = $true

$ws.Range("B3").Value = $true
$ws.Range("D3").Value = $true

$ws.Range("B4").Value = $false
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = $false

# Column widths (best-fit autofit widths from the source file; the
# runtime quantizes ColumnWidth to 1/6-character steps, so feed it the
# pre-image that lands on the nearest achievable value).
$ws.Range("B1").ColumnWidth = 17.498697916666668
$ws.Range("C1").ColumnWidth = 21.166666666666668
$ws.Range("D1").ColumnWidth = 19.330729166666668

# Selection
$ws.Range("C4").Select()
